# Update "想去人数" (F column) counts across sheets per gh-pages rebuild.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 174
$ws1.Range("F3").Value = 486
$ws1.Range("F4").Value = 18
$ws1.Range("F7").Value = 30
$ws1.Range("F8").Value = 19
$ws1.Range("F9").Value = 947

# Sheet "演出" (sheet2)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 88

# Sheet "全部类型" (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 174
$ws4.Range("F3").Value = 88
$ws4.Range("F4").Value = 486
$ws4.Range("F5").Value = 18
$ws4.Range("F8").Value = 30
$ws4.Range("F9").Value = 19
$ws4.Range("F10").Value = 947
